$d = $word.ActiveDocument

# --- Part 1: move the "_GoBack" bookmark away from the document start ---
# (it currently wraps the very beginning of "CONTRATO DE EMPLEO"; the edit
# relocates it into the last paragraph, between the new "01" and "sit").
# We'll re-add it later at the correct spot - re-adding a bookmark with the
# same name simply relocates it, so nothing else is required here.

# --- Part 2: rewrite the start of the last body paragraph ---
# "Lorem ipsum dolor sit amet..." -> "Modificada 01<GoBack/>sit amet..."
# with "01" and "sit" underlined.

$rng = $d.Content
$found = $rng.Find.Execute("Lorem ipsum dolor sit ")
if (-not $found) {
    Write-Host "ERROR: target text not found"
} else {
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # Delete "Lorem ipsum dolor " (keep the trailing "sit ")
    $prefix = $d.Range($matchStart, $matchEnd - 4)
    $prefix.Text = ""

    # Underline the kept "sit" (3 chars right after the deletion point)
    $sitWord = $d.Range($matchStart, $matchStart + 3)
    $sitWord.Font.Underline = 1

    # Insert the bookmark right before "sit" (zero-length range)
    $bm = $d.Range($matchStart, $matchStart)
    $d.Bookmarks.Add("_GoBack", $bm)

    # Insert "01" (underlined) right before the bookmark/"sit"
    $ins01 = $d.Range($matchStart, $matchStart)
    $ins01.InsertBefore("01")
    $ins01.Font.Underline = 1

    # Insert " " (not underlined) before "01"
    $insSpace = $d.Range($matchStart, $matchStart)
    $insSpace.InsertBefore(" ")

    # Insert "Modificada" (not underlined) before the space
    $insWord = $d.Range($matchStart, $matchStart)
    $insWord.InsertBefore("Modificada")

    Write-Host "rewrote paragraph start"
}
